$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.929.47"
$ws.Range("E2").Value = "  +1.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.297.27"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.54"
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.44"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +1.59%  "

$ws.Range("E9").Value = "  +4.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  -1.85%  "

$ws.Range("E11").Value = "  +1.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.874.12"
$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("E14").Value = "  +3.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.937.68"
$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("E16").Value = "  +2.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.279.97"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.65"
$ws.Range("E20").Value = "  +5.16%  "

$ws.Range("E21").Value = "  +1.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.81"
$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.518"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000122"
$ws.Range("E25").Value = "  +1.58%  "

$ws.Range("E26").Value = "  +4.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.71"
$ws.Range("E27").Value = "  +0.49%  "

$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("E29").Value = "  +0.62%  "

$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.15"
$ws.Range("E31").Value = "  +2.05%  "

$ws.Range("E32").Value = "  +2.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.18"
$ws.Range("E33").Value = "  +4.25%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  +2.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.86"
$ws.Range("E36").Value = "  +0.83%  "

$ws.Range("E37").Value = "  +1.54%  "

$ws.Range("E38").Value = "  -2.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.61"
$ws.Range("E39").Value = "  +3.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.46"
$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.61"
$ws.Range("E41").Value = "  -3.41%  "

$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.59"
$ws.Range("E43").Value = "  +2.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0692"
$ws.Range("E44").Value = "  +2.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "343.67"
$ws.Range("E45").Value = "  -5.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.21"
$ws.Range("E46").Value = "  -1.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.622.89"
$ws.Range("E47").Value = "  -4.72%  "

$ws.Range("E48").Value = "  +1.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.15"
$ws.Range("E49").Value = "  +3.99%  "

$ws.Range("E50").Value = "  +3.15%  "

$ws.Range("E51").Value = "  +0.21%  "
